$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "1 child" row (row 8) counts per the "1205 cap 3 sum" fix.
$ws.Range("B8").Value = 881
$ws.Range("C8").Value = 3196
$ws.Range("D8").Value = 2194
$ws.Range("E8").Value = 312
$ws.Range("F8").Value = 41

# Recalculate so the dependent totals (G8) and percentage formulas
# (B18:F18), which the chart is built on, pick up the new figures.
$excel.CalculateFullRebuild()

# Match the selection left active on the sheet after the edit.
$null = $ws.Range("B2:F8").Select()
